$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the A/B data columns (rows 2-11) with the new random values
$ws.Range("A2").Value = 59
$ws.Range("B2").Value = 55

$ws.Range("A3").Value = 41
$ws.Range("B3").Value = 29

$ws.Range("A4").Value = 5
$ws.Range("B4").Value = 86

$ws.Range("A5").Value = 47
$ws.Range("B5").Value = 24

$ws.Range("A6").Value = -3
$ws.Range("B6").Value = 16

$ws.Range("A7").Value = 73
$ws.Range("B7").Value = 10

$ws.Range("A8").Value = 31
$ws.Range("B8").Value = 19

$ws.Range("A9").Value = 58
$ws.Range("B9").Value = 81

$ws.Range("A10").Value = 54
$ws.Range("B10").Value = 72

$ws.Range("A11").Value = 56
$ws.Range("B11").Value = 8

# Update the view: scroll back to top-left A1 and select I6
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("I6").Select()
